$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$row = 56
$c3 = $ws.Cells.Item($row, 3)
$c3.NumberFormat = "@"

$ws.Cells.Item($row, 1).Value = "2024-09-25T18:06:40Z"
$ws.Cells.Item($row, 2).Value = "temperature"
$c3.Value = "25"
$ws.Cells.Item($row, 4).Value = "N/A"
$ws.Cells.Item($row, 5).Value = "N/A"
$ws.Cells.Item($row, 6).Value = "N/A"

$c3.Style = "Normal"
